$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.150.23"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "2.618.80"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'588.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").Value = "'165.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("D9").Value = "2.618.47"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  -3.76%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "'27.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.79%  "
$ws.Range("D15").Value = "3.097.16"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "'0.0000180"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("D17").Value = "67.098.65"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "2.620.70"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "'11.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("E20").Value = "  -6.28%  "
$ws.Range("D21").Value = "'356.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").Value = "'4.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.01%  "
$ws.Range("D23").Value = "'4.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").Value = "'10.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.40%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -4.72%  "
$ws.Range("D27").Value = "'69.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "0.0₃0997"
$ws.Range("D31").Value = "'545.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  -4.03%  "
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -4.36%  "
$ws.Range("D38").Value = "'158.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").Value = "'0.364"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("D41").Value = "'18.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").Value = "'5.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.13%  "
$ws.Range("D45").Value = "'2.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").Value = "'0.579"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("D48").Value = "'151.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("D49").Value = "'3.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.20%  "
$ws.Range("D50").Value = "'1.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("E51").Value = "  -1.76%  "
